$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.380.17"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "1.825.86"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4465"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3760"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8909"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.49%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "1.828.98"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.763"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.416"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07110"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008808"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "27.385.02"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.264"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").Value = "2.057.00"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.975"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.381"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.366"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08828"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7872"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.202"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.522"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.904"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.112"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01993"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05338"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.392"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5319"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1733"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.286"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.757"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5115"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.705"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06377"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.67%  "
